# v2.3 Change owner state for ADMIN HOME and PUBLISH VIDEO features
#
# This script applies the edits captured by the commit:
#  - REVIEW-SHEET: normalize Owner names that were entered as
#    shorthand/lowercase ("Omar" -> "omar sherif", "hala" -> "Hala Eldaly")
#    to match the full names used elsewhere in the sheet.
#  - REVIEW-SHEET: close out the Owner Status for the ADMIN HOME and
#    PUBLISH VIDEO review rows (rows 19-21).
#  - VERSION-HISTORY: capitalize the "Close owner state for USERHOME
#    feature" entry, renumber the trailing v1.10/v1.11/v1.12 entries to
#    v2.0/v2.1/v2.2, and append a new v2.3 entry documenting this change.

$wb = $excel.ActiveWorkbook

$reviewSheet = $wb.Worksheets.Item("REVIEW-SHEET")
$versionSheet = $wb.Worksheets.Item("VERSION-HISTORY")

# --- REVIEW-SHEET: Owner column name normalization -------------------------
$reviewSheet.Range("H2").Value = "omar sherif"
$reviewSheet.Range("H3").Value = "omar sherif"
$reviewSheet.Range("H4").Value = "Hala Eldaly"
$reviewSheet.Range("H5").Value = "Hala Eldaly"

# --- REVIEW-SHEET: Owner Status for ADMIN HOME / PUBLISH VIDEO rows --------
$reviewSheet.Range("I19").Value = "not applicable"
$reviewSheet.Range("I20").Value = "closed"
$reviewSheet.Range("I21").Value = "closed"

# --- VERSION-HISTORY: fix capitalization of v1.8 entry ---------------------
$versionSheet.Range("C10").Value = "Close owner state for USERHOME  feature "

# --- VERSION-HISTORY: renumber v1.10/v1.11/v1.12 to v2.0/v2.1/v2.2 ---------
$versionSheet.Range("A12").Value = "v2.0"
$versionSheet.Range("A13").Value = "v2.1"
$versionSheet.Range("A14").Value = "v2.2"

# --- VERSION-HISTORY: append new v2.3 row ----------------------------------
$versionSheet.Range("A14:D14").Copy()
$versionSheet.Range("A15:D15").PasteSpecial(-4122)
$versionSheet.Range("A15").Value = "v2.3"
$versionSheet.Range("B15").Value = "Hala Eldaly"
$versionSheet.Range("C15").Value = "Change owner state for ADMIN HOME and PUBLISH VIDEO features"
$versionSheet.Range("D15").Value = Get-Date -Year 2025 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0
$versionSheet.Rows.Item(15).RowHeight = 37.5
